$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2539.8
$ws.Range("I51").Value = 3749.5
$ws.Range("J51").Value = 1733.3334
$ws.Range("K51").Value = 3749.5
$ws.Range("L51").Value = 1733.3334
$ws.Range("M51").Value = -3265.5
$ws.Range("N51").Value = -2701.3334
$ws.Range("H96").Value = 41667680
$ws.Range("I96").Value = 50000904
$ws.Range("K96").Value = 150002712
$ws.Range("M96").Value = -150001339
$ws.Range("H116").Value = 3086.9443
$ws.Range("I116").Value = 1365
$ws.Range("J116").Value = 4182.727
$ws.Range("K116").Value = 1365
$ws.Range("L116").Value = 4182.727
$ws.Range("M116").Value = 2077
$ws.Range("N116").Value = -11066.727
$ws.Range("H132").Value = 4173.174
$ws.Range("I132").Value = 4553.9443
$ws.Range("K132").Value = 13661.8329
$ws.Range("M132").Value = -11131.8329
$ws.Range("H135").Value = 25010136
$ws.Range("I135").Value = 1246.1538
$ws.Range("K135").Value = 11215.3842
$ws.Range("M135").Value = -8680.3842
$ws.Range("H137").Value = 40746.19
$ws.Range("I137").Value = 2515.625
$ws.Range("J137").Value = 101915.1
$ws.Range("K137").Value = 7546.875
$ws.Range("L137").Value = 305745.3
$ws.Range("M137").Value = -4996.875
$ws.Range("N137").Value = -310845.3
$ws.Range("H141").Value = 3001.3635
$ws.Range("I141").Value = 2265.8333
$ws.Range("J141").Value = 3884
$ws.Range("K141").Value = 6797.499899999999
$ws.Range("L141").Value = 11652
$ws.Range("M141").Value = -1617.499899999999
$ws.Range("N141").Value = -22012

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23665.291
$ws.Range("I32").Value = 23423.83
$ws.Range("K32").Value = 23423.83
$ws.Range("M32").Value = -23136.83
$ws.Range("H45").Value = 2678.4443
$ws.Range("I45").Value = 3437.25
$ws.Range("J45").Value = 2071.4
$ws.Range("K45").Value = 3437.25
$ws.Range("L45").Value = 2071.4
$ws.Range("M45").Value = -3060.25
$ws.Range("N45").Value = -2825.4
$ws.Range("H61").Value = 2426.4375
$ws.Range("I61").Value = 1866.75
$ws.Range("J61").Value = 3359.25
$ws.Range("K61").Value = 1866.75
$ws.Range("L61").Value = 3359.25
$ws.Range("M61").Value = -1654.75
$ws.Range("N61").Value = -3783.25
$ws.Range("H92").Value = 26499.5
$ws.Range("J92").Value = 26499.5
$ws.Range("L92").Value = 26499.5
$ws.Range("N92").Value = -31491.5
$ws.Range("H97").Value = 949.0909
$ws.Range("I97").Value = 869.4706
$ws.Range("J97").Value = 1219.8
$ws.Range("K97").Value = 869.4706
$ws.Range("L97").Value = 1219.8
$ws.Range("M97").Value = -373.4706
$ws.Range("N97").Value = -2211.8
$ws.Range("H101").Value = 50000
$ws.Range("J101").Value = 50000
$ws.Range("L101").Value = 50000
$ws.Range("N101").Value = -56490
$ws.Range("H122").Value = 1635.8182
$ws.Range("I122").Value = 1699.4
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 5098.200000000001
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -2648.200000000001
$ws.Range("N122").Value = -7900
$ws.Range("H132").Value = 12635.957
$ws.Range("I132").Value = 1712.9412
$ws.Range("J132").Value = 41203.848
$ws.Range("K132").Value = 5138.8236
$ws.Range("L132").Value = 123611.544
$ws.Range("M132").Value = -2608.8236
$ws.Range("N132").Value = -128671.544
$ws.Range("H136").Value = 2426.4375
$ws.Range("I136").Value = 1866.75
$ws.Range("J136").Value = 3359.25
$ws.Range("K136").Value = 5600.25
$ws.Range("L136").Value = 10077.75
$ws.Range("M136").Value = -3050.25
$ws.Range("N136").Value = -15177.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H100").Value = 23262
$ws.Range("J100").Value = 23262
$ws.Range("L100").Value = 23262
$ws.Range("N100").Value = -25426
$ws.Range("H107").Value = 1088
$ws.Range("I107").Value = 1066.6666
$ws.Range("K107").Value = 1066.6666
$ws.Range("M107").Value = 853.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 325.77777
$ws.Range("I22").Value = 186
$ws.Range("K22").Value = 186
$ws.Range("M22").Value = 164
$ws.Range("H96").Value = 11624
$ws.Range("J96").Value = 11624
$ws.Range("L96").Value = 11624
$ws.Range("N96").Value = -17116
$ws.Range("H107").Value = 636.6667
$ws.Range("I107").Value = 536.6429000000001
$ws.Range("J107").Value = 744.38464
$ws.Range("K107").Value = 536.6429000000001
$ws.Range("L107").Value = 744.38464
$ws.Range("M107").Value = 1383.3571
$ws.Range("N107").Value = -4584.38464
$ws.Range("H122").Value = 1231.909
$ws.Range("I122").Value = 1410.4
$ws.Range("J122").Value = 1083.1666
$ws.Range("K122").Value = 4231.200000000001
$ws.Range("L122").Value = 3249.4998
$ws.Range("M122").Value = -1781.200000000001
$ws.Range("N122").Value = -8149.4998
$ws.Range("H132").Value = 33960
$ws.Range("I132").Value = 45009
$ws.Range("J132").Value = 7442.4
$ws.Range("K132").Value = 135027
$ws.Range("L132").Value = 22327.2
$ws.Range("M132").Value = -132497
$ws.Range("N132").Value = -27387.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 833.75
$ws.Range("I113").Value = 900
$ws.Range("J113").Value = 824.2857
$ws.Range("K113").Value = 2700
$ws.Range("L113").Value = 2472.8571
$ws.Range("M113").Value = -530
$ws.Range("N113").Value = -6812.8571
$ws.Range("H131").Value = 759.88
$ws.Range("J131").Value = 763.80414
$ws.Range("L131").Value = 2291.41242
$ws.Range("N131").Value = -12371.41242

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3611.2632
$ws.Range("I80").Value = 3316.6667
$ws.Range("J80").Value = 3747.2307
$ws.Range("K80").Value = 3316.6667
$ws.Range("L80").Value = 3747.2307
$ws.Range("M80").Value = -2318.6667
$ws.Range("N80").Value = -5743.2307
$ws.Range("H83").Value = 3611.2632
$ws.Range("I83").Value = 3316.6667
$ws.Range("J83").Value = 3747.2307
$ws.Range("K83").Value = 16583.3335
$ws.Range("L83").Value = 18736.1535
$ws.Range("M83").Value = -11591.3335
$ws.Range("N83").Value = -28720.1535
$ws.Range("H132").Value = 55462.656
$ws.Range("I132").Value = 63327.293
$ws.Range("J132").Value = 44321.082
$ws.Range("K132").Value = 189981.879
$ws.Range("L132").Value = 132963.246
$ws.Range("M132").Value = -187451.879
$ws.Range("N132").Value = -138023.246

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 63.57143
$ws.Range("I55").Value = 33.76923
$ws.Range("K55").Value = 33.76923
$ws.Range("M55").Value = 139.23077
$ws.Range("H132").Value = 2805.353
$ws.Range("I132").Value = 1844.3636
$ws.Range("J132").Value = 4567.1665
$ws.Range("K132").Value = 5533.0908
$ws.Range("L132").Value = 13701.4995
$ws.Range("M132").Value = -3003.0908
$ws.Range("N132").Value = -18761.4995
$ws.Range("H136").Value = 1937.5625
$ws.Range("I136").Value = 1300.1
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 3900.3
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -1350.3
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 17602
$ws.Range("J101").Value = 17602
$ws.Range("L101").Value = 17602
$ws.Range("N101").Value = -24092
$ws.Range("H123").Value = 42900
$ws.Range("J123").Value = 42900
$ws.Range("L123").Value = 42900
$ws.Range("N123").Value = -52700
$ws.Range("H132").Value = 2534.3635
$ws.Range("I132").Value = 1580
$ws.Range("J132").Value = 3329.6667
$ws.Range("K132").Value = 4740
$ws.Range("L132").Value = 9989.000100000001
$ws.Range("M132").Value = -2210
$ws.Range("N132").Value = -15049.0001
$ws.Range("H136").Value = 38463916
$ws.Range("I136").Value = 76925590
$ws.Range("J136").Value = 2239.1538
$ws.Range("K136").Value = 230776770
$ws.Range("L136").Value = 6717.4614
$ws.Range("M136").Value = -230774220
$ws.Range("N136").Value = -11817.4614
